$d = $word.ActiveDocument

# --- Change 1 (do first, while "50M" is not yet present anywhere else in
#     the document): Siege Analytics bullet point.
# Splits the single run into three runs so "50M" gets its own bold +
# colored (2C3E50) run, matching the styling used for the "23%" / "64%"
# figures elsewhere in the same bullet.
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial machine",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial machine", 2)

$r = $d.Content
$r.Find.Execute("50M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = 1
$r.Font.Color = 5258796   # RGB(0x2C, 0x3E, 0x50) == wdColor 2C3E50

# --- Change 2: Professional summary paragraph ---
# "...affecting all Black and Asian-American voters, developed geospatial ML..."
#   -> "...affecting 50M voters, developed geospatial ML..."
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2)

# --- Change 3: "Impact:" paragraph for the redistricting project ---
# "...Corrected demographic data affecting all Black and Asian-American voters, improved..."
#   -> "...Corrected demographic data affecting 50M voters nationwide, improved..."
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, improved electoral",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters nationwide, improved electoral", 2)
